# Updated cryptos list (mirrors the GitHub Actions scraper refresh):
# refreshed Price (col D) / Volume(1h) (col E) figures, and restored the
# correct Bittensor / Binance-Peg BSC-USD row ordering + values in rows 29-30.
#
# Price-column values are written with a leading literal apostrophe
# ("'591.65") -- Excel's standard "treat as text" quote-prefix -- so that
# numeric-looking prices (e.g. "591.65", "1.00") stay plain text instead of
# being auto-coerced to numbers, matching how the source sheet stores them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''67.320.20'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").Value = '''2.610.07'
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''591.65'
$ws.Range("E5").Value = '  -1.81%  '
$ws.Range("D6").Value = '''150.41'
$ws.Range("E6").Value = '  -2.63%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.50%  '
$ws.Range("D9").Value = '''2.607.62'
$ws.Range("E9").Value = '  -0.34%  '
$ws.Range("E10").Value = '  +1.01%  '
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("D12").Value = '''5.17'
$ws.Range("E12").Value = '  -1.53%  '
$ws.Range("E13").Value = '  -3.19%  '
$ws.Range("D14").Value = '''27.25'
$ws.Range("E14").Value = '  -2.72%  '
$ws.Range("D15").Value = '''3.082.64'
$ws.Range("E15").Value = '  -0.41%  '
$ws.Range("E16").Value = '  -2.79%  '
$ws.Range("D17").Value = '''67.526.12'
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("D18").Value = '''2.608.15'
$ws.Range("E18").Value = '  -0.54%  '
$ws.Range("D19").Value = '''373.12'
$ws.Range("E19").Value = '  +2.47%  '
$ws.Range("D20").Value = '''11.03'
$ws.Range("E20").Value = '  -2.41%  '
$ws.Range("D21").Value = '''7.37'
$ws.Range("E21").Value = '  -3.28%  '
$ws.Range("E22").Value = '  -0.57%  '
$ws.Range("D23").Value = '''4.83'
$ws.Range("E23").Value = '  -2.95%  '
$ws.Range("D24").Value = '''2.04'
$ws.Range("E24").Value = '  -4.06%  '
$ws.Range("D25").Value = '''73.74'
$ws.Range("E25").Value = '  +5.16%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").Value = '''9.93'
$ws.Range("E27").Value = '  -1.86%  '
$ws.Range("D28").Value = '''2.738.13'
$ws.Range("E28").Value = '  -0.26%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  -2.91%  '
$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D30").Value = '''577.94'
$ws.Range("E30").Value = '  -2.27%  '
$ws.Range("D31").Value = '''0.0₃0985'
$ws.Range("E31").Value = '  -6.42%  '
$ws.Range("E32").Value = '  -5.72%  '
$ws.Range("D33").Value = '''7.66'
$ws.Range("E33").Value = '  -3.60%  '
$ws.Range("E34").Value = '  -3.43%  '
$ws.Range("D36").Value = '''0.126'
$ws.Range("E36").Value = '  -3.70%  '
$ws.Range("E37").Value = '  -3.28%  '
$ws.Range("D38").Value = '''157.74'
$ws.Range("E38").Value = '  +0.80%  '
$ws.Range("E39").Value = '  -2.03%  '
$ws.Range("E40").Value = '  -1.89%  '
$ws.Range("E41").Value = '  -0.68%  '
$ws.Range("D42").Value = '''5.23'
$ws.Range("E42").Value = '  -3.85%  '
$ws.Range("D43").Value = '''2.55'
$ws.Range("E43").Value = '  -4.28%  '
$ws.Range("D44").Value = '''17.13'
$ws.Range("E44").Value = '  +4.26%  '
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").Value = '''153.01'
$ws.Range("E46").Value = '  -2.53%  '
$ws.Range("D47").Value = '''0.0₆0284'
$ws.Range("E47").Value = '  -1.57%  '
$ws.Range("E48").Value = '  -1.50%  '
$ws.Range("D49").Value = '''0.0777'
$ws.Range("E49").Value = '  -1.83%  '
$ws.Range("E50").Value = '  -4.95%  '
$ws.Range("D51").Value = '''21.34'
$ws.Range("E51").Value = '  +1.22%  '
